$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# --- Update Sheet1 view state (selection / scroll) ---
$sheet1.Range("A31").Select() | Out-Null

# --- Add Sheet2 after Sheet1 ---
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "Sheet2"

# --- Populate Sheet2 data table ---
$ws.Range('B1').Value = 'VEHouseholdVehicles'

$ws.Range('B2').Value = 'R script'
$ws.Range('C2').Value = 'Model'
$ws.Range('D2').Value = 'Outcome'
$ws.Range('E2').Value = 'Notes'

$ws.Range('B3').Value = 'AdjustVehicleOwnership'
$ws.Range('C3').Value = 'None (adjustment)'
$ws.Range('D3').Value = 'Success'

$ws.Range('B4').Value = 'AssignDrivers'
$ws.Range('C4').Value = 'lm in estimateDriverModel'
$ws.Range('D4').Value = 'Warning message: In predict.lm(object, newdata, se.fit, scale = 1, type = if (type ==  :  prediction from a rank-deficient fit may be isleading'
$ws.Range('E4').Value = 'Missing Hometype causes error.'

$ws.Range('B5').Value = 'AssignVehicleAge'
$ws.Range('C5').Value = 'None (assignment)'
$ws.Range('D5').Value = 'Success with minor modification to "vehyear"'

$ws.Range('B6').Value = 'AssignVehicleFeatures'
$ws.Range('C6').Value = 'None (assignment)'
$ws.Range('D6').Value = 'Success'

$ws.Range('B7').Value = 'AssignVehicleFeaturesFuture'
$ws.Range('C7').Value = 'None (assignment)'
$ws.Range('D7').Value = 'Success'

$ws.Range('B8').Value = 'AssignVehicleOwnership'
$ws.Range('C8').Value = 'clm & glm'
$ws.Range('D8').Value = 'Success with warning on AutoOwnModels_ls$Metro$Zero <- glm()'
$ws.Range('E8').Value = 'Missing Hometype causes error'

$ws.Range('B9').Value = 'AssignVehicleType'
$ws.Range('C9').Value = 'glm estimateVehicleTypeModel'
$ws.Range('D9').Value = 'Error (Error in ''binarySearch'' function to match target value)'
$ws.Range('E9').Value = 'Missing Hometype coefficient'

$ws.Range('B10').Value = 'CalculateVehicleOwnCost'
$ws.Range('C10').Value = 'None (calculation)'
$ws.Range('D10').Value = 'Success'

$ws.Range('B11').Value = 'CreateVehicleTable'
$ws.Range('C11').Value = 'None (assignment)'
$ws.Range('D11').Value = 'Success'

$ws.Range('B12').Value = 'Finalize'
$ws.Range('C12').Value = 'None (assignment)'
$ws.Range('D12').Value = 'Success'

$ws.Range('B14').Value = 'VEHouseholdTravel'

$ws.Range('B15').Value = 'R script'
$ws.Range('C15').Value = 'Model'
$ws.Range('D15').Value = 'Outcome'
$ws.Range('E15').Value = 'Notes'

$ws.Range('B16').Value = 'ApplyDvmtReductions'
$ws.Range('C16').Value = 'None (assignment)'
$ws.Range('D16').Value = 'Success'

$ws.Range('B17').Value = 'CalculateAltModeTrips'
$ws.Range('C17').Value = 'Hurdle'
$ws.Range('D17').Value = 'Success'

$ws.Range('B18').Value = 'CalculateHouseholdDvmt'
$ws.Range('C18').Value = 'BinarySearch'
$ws.Range('D18').Value = 'Success'

$ws.Range('B19').Value = 'CalculateInducedDemand'
$ws.Range('C19').Value = 'None (calculation)'
$ws.Range('D19').Value = 'Success'

$ws.Range('B20').Value = 'CalculatePolicyVmt'
$ws.Range('C20').Value = 'None (calculation)'
$ws.Range('D20').Value = 'Success'

$ws.Range('B21').Value = 'CalculateTravelDemand'
$ws.Range('C21').Value = 'None (calculation)'
$ws.Range('D21').Value = 'Success'

$ws.Range('B22').Value = 'CalculateTravelDemandFuture'
$ws.Range('C22').Value = 'None (calculation)'
$ws.Range('D22').Value = 'Success'

$ws.Range('B23').Value = 'CalculateVehicleTrips'
$ws.Range('C23').Value = 'None (calculation)'
$ws.Range('D23').Value = 'Success'

$ws.Range('B24').Value = 'DivertSovTravel'
$ws.Range('C24').Value = 'ZeroSov_GLM'
$ws.Range('D24').Value = 'Success'

$ws.Range('B25').Value = 'Initialize'
$ws.Range('C25').Value = 'None'
$ws.Range('D25').Value = 'Success'

$ws.Range('B26').Value = 'LoadDefaultValues'
$ws.Range('C26').Value = 'None'
$ws.Range('D26').Value = 'Success'

$ws.Range('B28').Value = 'VEPowertrainsAndFuels'

$ws.Range('B29').Value = 'R script'
$ws.Range('C29').Value = 'Model'
$ws.Range('D29').Value = 'Outcome'
$ws.Range('E29').Value = 'Notes'

$ws.Range('B30').Value = 'Initialize'
$ws.Range('C30').Value = 'None'
$ws.Range('D30').Value = 'Success'

$ws.Range('B31').Value = 'LoadDefaultValues'
$ws.Range('C31').Value = 'None'
$ws.Range('D31').Value = 'Success'

$ws.Range('B32').Value = 'CalculateCarbonIntensity'
$ws.Range('C32').Value = 'None (calculation)'
$ws.Range('D32').Value = 'Success'

$ws.Range('B33').Value = 'AssignHhVehiclePowertrain'
$ws.Range('C33').Value = 'None (assignment)'
$ws.Range('D33').Value = 'Success'

$ws.Range('B35').Value = 'VELandUse'

$ws.Range('B36').Value = 'R script'
$ws.Range('C36').Value = 'Model'
$ws.Range('D36').Value = 'Outcome'
$ws.Range('E36').Value = 'Notes'

$ws.Range('B37').Value = 'AssignCarSvcAvailability'
$ws.Range('C37').Value = 'None (assignment)'
$ws.Range('D37').Value = 'Success'

$ws.Range('B38').Value = 'AssignDemandManagement'
$ws.Range('C38').Value = 'None (assignment)'
$ws.Range('D38').Value = 'Success'

$ws.Range('B39').Value = 'AssignDevTypes'
$ws.Range('C39').Value = 'None (assignment)'
$ws.Range('D39').Value = 'Success'

$ws.Range('B40').Value = 'AssignLocTypes'
$ws.Range('C40').Value = 'None (assignment)'
$ws.Range('D40').Value = 'Success'

$ws.Range('B41').Value = 'AssignParkingRestrictions'
$ws.Range('C41').Value = 'None (assignment)'
$ws.Range('D41').Value = 'Success'

$ws.Range('B42').Value = 'Calculate4DMeasures'
$ws.Range('C42').Value = 'None (calculation)'
$ws.Range('D42').Value = 'Success'

$ws.Range('B43').Value = 'CalculateBasePlaceTypes'
$ws.Range('C43').Value = 'None (calculation)'
$ws.Range('D43').Value = 'Success'

$ws.Range('B44').Value = 'CalculateFuturePlaceTypes'
$ws.Range('C44').Value = 'None (calculation)'
$ws.Range('D44').Value = 'Success'

$ws.Range('B45').Value = 'CalculateUrbanMixMeasure'
$ws.Range('C45').Value = 'CalculateUrbanMixMeasure binarySearch'
$ws.Range('D45').Value = 'Error (Error in ''binarySearch'' function to match target value)'
$ws.Range('E45').Value = 'Likely due to missing Hometype'

$ws.Range('B46').Value = 'LocateEmployment'
$ws.Range('C46').Value = 'None (assignment)'
$ws.Range('D46').Value = 'Success'

$ws.Range('B47').Value = 'PredictHousing'
$ws.Range('C47').Value = 'HouseTypeModel glm'
$ws.Range('D47').Value = 'Success'
$ws.Range('E47').Value = 'but doubt results given missing Hometype variable'

# --- Column widths (best effort match to target character widths) ---
$ws.Columns.Item(2).ColumnWidth = 27.666666666666668
$ws.Columns.Item(3).ColumnWidth = 37.5
$ws.Columns.Item(4).ColumnWidth = 45.0
$ws.Columns.Item(5).ColumnWidth = 13.333333333333334

# --- Selection on Sheet2 (becomes the active/visible sheet) ---
$ws.Range("B36:E36").Select() | Out-Null
